$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) "Converting blank values to NaN." - collapse the 3 runs (with
#    proofErr spell-check wrapper around "NaN") into a single run.
# ---------------------------------------------------------------
$p = $d.Paragraphs(60)
$r = $p.Range
$full = $d.Range($r.Start, $r.End - 1)
$full.Text = "TEMP1"
$full2 = $d.Range($r.Start, $r.Start + 5)
$full2.Text = "Converting blank values to NaN."

# ---------------------------------------------------------------
# 2) "Remove instances (rows) with NaN." - same simplification.
# ---------------------------------------------------------------
$p = $d.Paragraphs(61)
$r = $p.Range
$full = $d.Range($r.Start, $r.End - 1)
$full.Text = "TEMP2"
$full2 = $d.Range($r.Start, $r.Start + 5)
$full2.Text = "Remove instances (rows) with NaN."

# ---------------------------------------------------------------
# 3) "Use a classification Algorithm." -> "Use a Regression Algorithm."
#    split into 3 runs: "Use a " / "Regression" / " Algorithm."
# ---------------------------------------------------------------
$p = $d.Paragraphs(64)
$r = $p.Range
$start = $r.Start
$mid = $d.Range($start + 6, $start + 20)   # "classification"
$mid.Text = "Regression"
# Force a run boundary around the replaced word by toggling Bold off/on.
$mid2 = $d.Range($start + 6, $start + 16)  # "Regression"
$mid2.Bold = 1
$mid2.Bold = 0

# ---------------------------------------------------------------
# 4) "Random Forest classifier" -> "Linear Regression " split into
#    3 runs: "L" / "inear" / " Regression "
# ---------------------------------------------------------------
$p = $d.Paragraphs(65)
$r = $p.Range
$start = $r.Start
$end = $r.End
$full = $d.Range($start, $end - 1)
$full.Text = "Linear Regression "
# Force run boundaries: "L" | "inear" | " Regression "
$seg1 = $d.Range($start, $start + 1)        # "L"
$seg1.Bold = 1
$seg1.Bold = 0
$seg2 = $d.Range($start + 1, $start + 6)    # "inear"
$seg2.Bold = 1
$seg2.Bold = 0

# ---------------------------------------------------------------
# 5) Delete the empty bullet paragraph that followed
#    "Random Forest classifier" / now "Linear Regression ".
# ---------------------------------------------------------------
$empty = $d.Paragraphs(66)
$empty.Range.Delete()

# ---------------------------------------------------------------
# 6) The paragraph that used to read "It combines multiple decision
#    trees..." drops its numbered-list formatting and becomes an
#    indented plain paragraph instead.
# ---------------------------------------------------------------
$it = $d.Paragraphs(66)
$it.Range.ListFormat.RemoveNumbers()
$it.LeftIndent = 36
